$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. New "Revenue" column block (AD:AK) with a merged header and year
#    columns (matching the existing EV/EPS & EPS blocks in Q:V / W:AB).
#    Written first so the new shared string ("Revenue") lands at the same
#    index the author's workbook used.
# ---------------------------------------------------------------------------
$ws.Range("AD2:AK2").Merge()
$ws.Range("AD2").Value = "Revenue"
$ws.Range("AD2").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 2. Insert a new row for "Evolution Gaming" / "EVO ST" above the old row 5
#    (Flutter), pushing Flutter and everything below it down by one row.
# ---------------------------------------------------------------------------
$ws.Range("A5").EntireRow.Insert()

$ws.Range("C5").Value = "Evolution Gaming"
$ws.Range("D5").Value = "EVO ST"

# Repair the running index column (B) so the numbering is contiguous again
# after the insert (B5 keeps counting on from B4, etc.).
$ws.Cells.Item(5, 2).Formula = "=+B4+1"
for ($r = 6; $r -le 18; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 2).Formula = "=+B$prev+1"
}

$ws.Range("AD3").Value = 2019
$ws.Range("AE3").Value = 2020
$ws.Range("AF3").Value = 2021
$ws.Range("AD3:AF3").HorizontalAlignment = -4108   # xlCenter

$ws.Range("AG3").Formula = "=+AF3+1"
$ws.Range("AH3:AK3").Formula = "=+AG3+1"

# Revenue figures (Flutter, now on row 6)
$ws.Range("AD6").Value = 2140
$ws.Range("AE6").Value = 4414
$ws.Range("AF6").Value = 6036

# ---------------------------------------------------------------------------
# 3. View state: freeze panes at column E / row 4, selections as left by
#    the author.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("E4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F5").Select()
